$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '332.38'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '1.20%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '44.22'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '7.09%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.865'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '4.40%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08336'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '2.03%'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '0.83%'
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.973'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-2.45%'
$ws.Range('B8').Value = 'BTSEToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.898'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-2.19%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.9367'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '1.65%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1255'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-1.82%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1962'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '0.23%'
$ws.Range('B12').Value = 'MCDex'
$ws.Range('C12').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.945'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '7.40%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09654'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '3.10%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.03998'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.1066'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.65%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001305'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '0.39%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.006003'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-2.13%'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.505'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '1.82%'
$ws.Range('B19').Value = 'GateToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.503'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-0.40%'
$ws.Range('B20').Value = 'BitpandaEcosystemToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.3510'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '0.38%'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-0.58%'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '6.60%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04402'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-0.29%'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.09%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004436'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '3.10%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001191'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '0.74%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02814'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '1.18%'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05692'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '5.30%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007935'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '3.71%'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1428'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '0.88%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.009008'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '0.25%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002102'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-1.97%'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-9.52%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00007249'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '9.55%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.11%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.003236'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '1.11%'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.14%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002102'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.11%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002002'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.11%'
